# Added Test data for Spain Zettler Market
#
# Mirrors the "Italy" tab into a new "Spain" tab (Excel's usual
# copy-the-most-similar-country-sheet workflow), then swaps in the
# Spain-specific market name / part number and refreshes the view state
# so "Spain" becomes the active/selected sheet.

$wb = $excel.ActiveWorkbook

$italy = $wb.Worksheets.Item("Italy")

# Select the whole used range on Italy before copying away from it -
# this is what the workbook looks like afterwards (Italy is no longer
# the active tab, so its lingering single-cell selection is replaced by
# the full-range selection left over from building the new sheet).
$italy.Activate()
$italy.Range("A1:D13").Select() | Out-Null

# Duplicate "Italy" and place the copy right after it.
$italy.Copy($null, $italy)
$spain = $wb.Worksheets.Item($italy.Index + 1)
$spain.Name = "Spain"

# Spain-specific content.
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2064/T2063"

# The shorter market name/part number let column B narrow down (and
# column D widen slightly); with the narrower column B the wrapped
# description cells in rows 3-5 now take two lines.
$spain.Columns.Item(2).ColumnWidth = 20.1
$spain.Columns.Item(3).ColumnWidth = $spain.StandardWidth
$spain.Columns.Item(4).ColumnWidth = 23.15

$spain.Rows.Item(3).RowHeight = 28.8
$spain.Rows.Item(4).RowHeight = 28.8
$spain.Rows.Item(5).RowHeight = 28.8
$spain.Rows.Item(13).AutoFit()

# Make Spain the active sheet/tab with its own selection.
$spain.Activate()
$spain.Range("F5").Select() | Out-Null
